$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("E19").Value = "test2"
$ws.Range("E19").Select()
